# Scheduled-runner refresh of cached Universalis market data
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns H:N)
# across the per-job Sheets. Plain cached numbers, no formulas involved.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 392.93332
$ws.Range("I53").Value = 246
$ws.Range("J53").Value = 539.86664
$ws.Range("K53").Value = 246
$ws.Range("L53").Value = 539.86664
$ws.Range("M53").Value = 391
$ws.Range("H107").Value = 1545.4286
$ws.Range("I107").Value = 1721.2858
$ws.Range("J107").Value = 1193.7142
$ws.Range("K107").Value = 1721.2858
$ws.Range("L107").Value = 1193.7142
$ws.Range("M107").Value = 198.7141999999999
$ws.Range("N107").Value = -5033.7142
$ws.Range("H113").Value = 3535.476
$ws.Range("I113").Value = 1616.3334
$ws.Range("J113").Value = 4303.1333
$ws.Range("K113").Value = 1616.3334
$ws.Range("L113").Value = 4303.1333
$ws.Range("M113").Value = 1637.6666
$ws.Range("N113").Value = -10811.1333
$ws.Range("H132").Value = 17896244
$ws.Range("I132").Value = 24048580
$ws.Range("J132").Value = 669701.1
$ws.Range("K132").Value = 72145740
$ws.Range("L132").Value = 2009103.3
$ws.Range("M132").Value = -72143210
$ws.Range("N132").Value = -2014163.3
$ws.Range("H137").Value = 2365.25
$ws.Range("I137").Value = 1106.8286
$ws.Range("J137").Value = 5753.3076
$ws.Range("K137").Value = 3320.4858
$ws.Range("L137").Value = 17259.9228
$ws.Range("M137").Value = -770.4858000000004
$ws.Range("N137").Value = -22359.9228
$ws.Range("H138").Value = 3125.3225
$ws.Range("I138").Value = 768.9375
$ws.Range("J138").Value = 4361.459
$ws.Range("K138").Value = 2306.8125
$ws.Range("L138").Value = 13084.377
$ws.Range("M138").Value = 2833.1875
$ws.Range("N138").Value = -23364.377
$ws.Range("H141").Value = 6464.8687
$ws.Range("I141").Value = 6509.8647
$ws.Range("J141").Value = 4800
$ws.Range("K141").Value = 19529.5941
$ws.Range("L141").Value = 14400
$ws.Range("M141").Value = -14349.5941
$ws.Range("N141").Value = -24760

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1677.6666
$ws.Range("I2").Value = 1637.375
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 1637.375
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -1524.375
$ws.Range("N2").Value = -2226
$ws.Range("H32").Value = 5613.705
$ws.Range("I32").Value = 4408.1704
$ws.Range("J32").Value = 9660.857
$ws.Range("K32").Value = 4408.1704
$ws.Range("L32").Value = 9660.857
$ws.Range("M32").Value = -4121.1704
$ws.Range("N32").Value = -10234.857
$ws.Range("H74").Value = 2272.8909
$ws.Range("I74").Value = 2069.2654
$ws.Range("J74").Value = 3935.8333
$ws.Range("K74").Value = 2069.2654
$ws.Range("L74").Value = 3935.8333
$ws.Range("M74").Value = -1195.2654
$ws.Range("N74").Value = -5683.8333
$ws.Range("H77").Value = 2272.8909
$ws.Range("I77").Value = 2069.2654
$ws.Range("J77").Value = 3935.8333
$ws.Range("K77").Value = 10346.327
$ws.Range("L77").Value = 19679.1665
$ws.Range("M77").Value = -5978.327000000001
$ws.Range("N77").Value = -28415.1665
$ws.Range("H110").Value = 1929.3077
$ws.Range("I110").Value = 1944.375
$ws.Range("J110").Value = 1905.2
$ws.Range("K110").Value = 1944.375
$ws.Range("L110").Value = 1905.2
$ws.Range("M110").Value = 100.625
$ws.Range("N110").Value = -5995.2
$ws.Range("H116").Value = 1677.6666
$ws.Range("I116").Value = 1637.375
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 1637.375
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = 656.625
$ws.Range("N116").Value = -6588
$ws.Range("H132").Value = 1833.7954
$ws.Range("I132").Value = 1134
$ws.Range("J132").Value = 4555.222
$ws.Range("K132").Value = 3402
$ws.Range("L132").Value = 13665.666
$ws.Range("M132").Value = -872
$ws.Range("N132").Value = -18725.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1677.6666
$ws.Range("I3").Value = 1637.375
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 1637.375
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -1523.375
$ws.Range("N3").Value = -2228

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8335229.5
$ws.Range("I31").Value = 1044.619
$ws.Range("J31").Value = 27781660
$ws.Range("K31").Value = 1044.619
$ws.Range("L31").Value = 27781660
$ws.Range("M31").Value = -749.6189999999999
$ws.Range("N31").Value = -27782250
$ws.Range("H34").Value = 8335229.5
$ws.Range("I34").Value = 1044.619
$ws.Range("J34").Value = 27781660
$ws.Range("K34").Value = 1044.619
$ws.Range("L34").Value = 27781660
$ws.Range("M34").Value = -842.6189999999999
$ws.Range("N34").Value = -27782064
$ws.Range("H99").Value = 8701093
$ws.Range("I99").Value = 18185512
$ws.Range("J99").Value = 7041.6665
$ws.Range("K99").Value = 18185512
$ws.Range("L99").Value = 7041.6665
$ws.Range("M99").Value = -18184014
$ws.Range("N99").Value = -10037.6665
$ws.Range("H126").Value = 8701093
$ws.Range("I126").Value = 18185512
$ws.Range("J126").Value = 7041.6665
$ws.Range("K126").Value = 54556536
$ws.Range("L126").Value = 21124.9995
$ws.Range("M126").Value = -54554066
$ws.Range("N126").Value = -26064.9995
$ws.Range("H132").Value = 1834.9166
$ws.Range("I132").Value = 971.5185
$ws.Range("J132").Value = 4425.1113
$ws.Range("K132").Value = 2914.5555
$ws.Range("L132").Value = 13275.3339
$ws.Range("M132").Value = -384.5554999999999
$ws.Range("N132").Value = -18335.3339
$ws.Range("H141").Value = 34933.332
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 34933.332
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 34933.332
$ws.Range("N141").Value = -45293.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 142.57895
$ws.Range("I33").Value = 149.375
$ws.Range("J33").Value = 137.63637
$ws.Range("K33").Value = 896.25
$ws.Range("L33").Value = 825.81822
$ws.Range("M33").Value = -613.25
$ws.Range("N33").Value = -1391.81822
$ws.Range("H113").Value = 678
$ws.Range("I113").Value = 724
$ws.Range("J113").Value = 650.4
$ws.Range("K113").Value = 2172
$ws.Range("L113").Value = 1951.2
$ws.Range("M113").Value = -2
$ws.Range("N113").Value = -6291.2
$ws.Range("H122").Value = 2809.5264
$ws.Range("I122").Value = 563.8333
$ws.Range("J122").Value = 3846
$ws.Range("K122").Value = 5074.4997
$ws.Range("L122").Value = 34614
$ws.Range("M122").Value = -2624.4997
$ws.Range("N122").Value = -39514

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 13621.223
$ws.Range("I113").Value = 1447.5
$ws.Range("J113").Value = 23360.2
$ws.Range("K113").Value = 1447.5
$ws.Range("L113").Value = 23360.2
$ws.Range("M113").Value = 722.5
$ws.Range("N113").Value = -27700.2
$ws.Range("H132").Value = 2350.851
$ws.Range("I132").Value = 1467.625
$ws.Range("J132").Value = 4235.067
$ws.Range("K132").Value = 4402.875
$ws.Range("L132").Value = 12705.201
$ws.Range("M132").Value = -1872.875
$ws.Range("N132").Value = -17765.201

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 825.5
$ws.Range("I16").Value = 825.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 825.5
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -655.5
$ws.Range("N16").ClearContents()
$ws.Range("H46").Value = 2850
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2850
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 2850
$ws.Range("N46").Value = -3226
$ws.Range("M46").ClearContents()
$ws.Range("H132").Value = 9982.697
$ws.Range("I132").Value = 11900.565
$ws.Range("J132").Value = 7777.15
$ws.Range("K132").Value = 35701.695
$ws.Range("L132").Value = 23331.45
$ws.Range("M132").Value = -33171.695
$ws.Range("N132").Value = -28391.45
$ws.Range("H136").Value = 2300.8
$ws.Range("I136").Value = 1168.96
$ws.Range("J136").Value = 7960
$ws.Range("K136").Value = 3506.88
$ws.Range("L136").Value = 23880
$ws.Range("M136").Value = -956.8800000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 563013.5
$ws.Range("I126").Value = 2328.3
$ws.Range("J126").Value = 1185997.1
$ws.Range("K126").Value = 6984.900000000001
$ws.Range("L126").Value = 3557991.3
$ws.Range("M126").Value = -4514.900000000001
